$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8407.277
$ws.Range("I86").Value = 2876
$ws.Range("K86").Value = 2876
$ws.Range("M86").Value = -1753
$ws.Range("H89").Value = 8407.277
$ws.Range("I89").Value = 2876
$ws.Range("K89").Value = 14380
$ws.Range("M89").Value = -8764
$ws.Range("H99").Value = 188
$ws.Range("I99").Value = 188
$ws.Range("K99").Value = 564
$ws.Range("M99").Value = 934
$ws.Range("H101").Value = 4092.5
$ws.Range("J101").Value = 4092.5
$ws.Range("L101").Value = 12277.5
$ws.Range("N101").Value = -15521.5
$ws.Range("H111").Value = 2351.5881
$ws.Range("I111").Value = 2676.9285
$ws.Range("J111").Value = 833.3333
$ws.Range("K111").Value = 8030.7855
$ws.Range("L111").Value = 2499.9999
$ws.Range("M111").Value = -4963.7855
$ws.Range("N111").Value = -8633.999899999999
$ws.Range("H112").Value = 1074.1818
$ws.Range("J112").Value = 1074.1818
$ws.Range("L112").Value = 3222.5454
$ws.Range("N112").Value = -5438.5454
$ws.Range("H115").Value = 612.5
$ws.Range("I115").Value = 612.5
$ws.Range("K115").Value = 1837.5
$ws.Range("M115").Value = -270.5
$ws.Range("H118").Value = 1067.5
$ws.Range("I118").Value = 726.6667
$ws.Range("K118").Value = 2180.0001
$ws.Range("M118").Value = -523.0001000000002
$ws.Range("H127").Value = 1276
$ws.Range("I127").Value = 881.8333
$ws.Range("J127").Value = 1613.8572
$ws.Range("K127").Value = 2645.4999
$ws.Range("L127").Value = 4841.571599999999
$ws.Range("M127").Value = 2314.5001
$ws.Range("N127").Value = -14761.5716
$ws.Range("H132").Value = 1906.8914
$ws.Range("I132").Value = 2028.7675
$ws.Range("K132").Value = 6086.3025
$ws.Range("M132").Value = -3556.3025
$ws.Range("H135").Value = 20006226
$ws.Range("I135").Value = 607.8
$ws.Range("K135").Value = 5470.2
$ws.Range("M135").Value = -2935.2
$ws.Range("H138").Value = 2343.0984
$ws.Range("I138").Value = 2579
$ws.Range("J138").Value = 2291.2
$ws.Range("K138").Value = 7737
$ws.Range("L138").Value = 6873.599999999999
$ws.Range("M138").Value = -2597
$ws.Range("N138").Value = -17153.6
$ws.Range("H141").Value = 1166.0476
$ws.Range("I141").Value = 910.4054
$ws.Range("K141").Value = 2731.2162
$ws.Range("M141").Value = 2448.7838

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6500
$ws.Range("I28").Value = 6500
$ws.Range("K28").Value = 6500
$ws.Range("M28").Value = -6308
$ws.Range("H32").Value = 5827.6924
$ws.Range("I32").Value = 6019.7676
$ws.Range("J32").Value = 4910
$ws.Range("K32").Value = 6019.7676
$ws.Range("L32").Value = 4910
$ws.Range("M32").Value = -5732.7676
$ws.Range("N32").Value = -5484
$ws.Range("H45").Value = 2643.516
$ws.Range("I45").Value = 2446.348
$ws.Range("J45").Value = 3210.375
$ws.Range("K45").Value = 2446.348
$ws.Range("L45").Value = 3210.375
$ws.Range("M45").Value = -2069.348
$ws.Range("N45").Value = -3964.375
$ws.Range("H99").Value = 6500
$ws.Range("I99").Value = 6500
$ws.Range("K99").Value = 6500
$ws.Range("M99").Value = -3505
$ws.Range("H110").Value = 609.8333
$ws.Range("I110").Value = 611.8
$ws.Range("K110").Value = 611.8
$ws.Range("M110").Value = 1433.2
$ws.Range("H132").Value = 27508.242
$ws.Range("I132").Value = 1379.5
$ws.Range("J132").Value = 127668.414
$ws.Range("K132").Value = 4138.5
$ws.Range("L132").Value = 383005.242
$ws.Range("M132").Value = -1608.5
$ws.Range("N132").Value = -388065.242

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3921.6365
$ws.Range("I105").Value = 4262.5713
$ws.Range("K105").Value = 4262.5713
$ws.Range("M105").Value = -2515.5713
$ws.Range("H134").Value = 7140.0625
$ws.Range("I134").Value = 8434.5
$ws.Range("J134").Value = 3256.75
$ws.Range("K134").Value = 25303.5
$ws.Range("L134").Value = 9770.25
$ws.Range("M134").Value = -22768.5
$ws.Range("N134").Value = -14840.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 85857290
$ws.Range("I6").Value = 15250007
$ws.Range("J6").Value = 180000340
$ws.Range("K6").Value = 15250007
$ws.Range("L6").Value = 180000340
$ws.Range("M6").Value = -15249894
$ws.Range("N6").Value = -180000566
$ws.Range("H31").Value = 12782.605
$ws.Range("I31").Value = 19321.545
$ws.Range("K31").Value = 19321.545
$ws.Range("M31").Value = -19026.545
$ws.Range("H34").Value = 12782.605
$ws.Range("I34").Value = 19321.545
$ws.Range("K34").Value = 19321.545
$ws.Range("M34").Value = -19119.545
$ws.Range("H122").Value = 1062.7878
$ws.Range("I122").Value = 920.8946999999999
$ws.Range("J122").Value = 1255.3572
$ws.Range("K122").Value = 2762.6841
$ws.Range("L122").Value = 3766.0716
$ws.Range("M122").Value = -312.6840999999999
$ws.Range("N122").Value = -8666.071599999999
$ws.Range("H132").Value = 17604.94
$ws.Range("I132").Value = 18825.414
$ws.Range("J132").Value = 8756.5
$ws.Range("K132").Value = 56476.242
$ws.Range("L132").Value = 26269.5
$ws.Range("M132").Value = -53946.242
$ws.Range("N132").Value = -31329.5
$ws.Range("H134").Value = 588.1111
$ws.Range("I134").Value = 486.69446
$ws.Range("K134").Value = 1460.08338
$ws.Range("M134").Value = 1074.91662

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 569.3333
$ws.Range("J52").Value = 569.3333
$ws.Range("L52").Value = 1707.9999
$ws.Range("N52").Value = -2239.9999
$ws.Range("H122").Value = 757.4783
$ws.Range("I122").Value = 357.5
$ws.Range("J122").Value = 898.64703
$ws.Range("K122").Value = 3217.5
$ws.Range("L122").Value = 8087.82327
$ws.Range("M122").Value = -767.5
$ws.Range("N122").Value = -12987.82327
$ws.Range("H129").Value = 385458.38
$ws.Range("J129").Value = 625874.9
$ws.Range("L129").Value = 1877624.7
$ws.Range("N129").Value = -1887624.7
$ws.Range("H131").Value = 119863.68
$ws.Range("J131").Value = 127419.23
$ws.Range("L131").Value = 382257.69
$ws.Range("N131").Value = -392337.69
$ws.Range("H137").Value = 41669844
$ws.Range("I137").Value = 1344.75
$ws.Range("J137").Value = 83338344
$ws.Range("K137").Value = 4034.25
$ws.Range("L137").Value = 250015032
$ws.Range("M137").Value = 1065.75
$ws.Range("N137").Value = -250025232
$ws.Range("H139").Value = 1242.9678
$ws.Range("I139").Value = 1165.0667
$ws.Range("J139").Value = 3580
$ws.Range("K139").Value = 3495.2001
$ws.Range("L139").Value = 10740
$ws.Range("M139").Value = 1644.7999
$ws.Range("N139").Value = -21020
$ws.Range("H140").Value = 2820.2727
$ws.Range("I140").Value = 1904.2858
$ws.Range("J140").Value = 4423.25
$ws.Range("K140").Value = 5712.857400000001
$ws.Range("L140").Value = 13269.75
$ws.Range("M140").Value = -532.8574000000008
$ws.Range("N140").Value = -23629.75
$ws.Range("H141").Value = 3814.75
$ws.Range("I141").Value = 3577.9
$ws.Range("K141").Value = 10733.7
$ws.Range("M141").Value = -5553.700000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9472.5
$ws.Range("I70").Value = 13291.8
$ws.Range("J70").Value = 4698.375
$ws.Range("K70").Value = 13291.8
$ws.Range("L70").Value = 4698.375
$ws.Range("M70").Value = -13021.8
$ws.Range("N70").Value = -5238.375
$ws.Range("H73").Value = 9472.5
$ws.Range("I73").Value = 13291.8
$ws.Range("J73").Value = 4698.375
$ws.Range("K73").Value = 13291.8
$ws.Range("L73").Value = 4698.375
$ws.Range("M73").Value = -12355.8
$ws.Range("N73").Value = -6570.375
$ws.Range("H102").Value = 45458456
$ws.Range("I102").Value = 45458456
$ws.Range("K102").Value = 45458456
$ws.Range("M102").Value = -45456834
$ws.Range("H122").Value = 66667836
$ws.Range("I122").Value = 25641712
$ws.Range("K122").Value = 76925136
$ws.Range("M122").Value = -76922686
$ws.Range("H132").Value = 15906.237
$ws.Range("I132").Value = 3001.1538
$ws.Range("J132").Value = 43867.25
$ws.Range("K132").Value = 9003.4614
$ws.Range("L132").Value = 131601.75
$ws.Range("M132").Value = -6473.4614
$ws.Range("N132").Value = -136661.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5884.143
$ws.Range("I22").Value = 5198.1665
$ws.Range("K22").Value = 5198.1665
$ws.Range("M22").Value = -4903.1665
$ws.Range("H27").Value = 5884.143
$ws.Range("I27").Value = 5198.1665
$ws.Range("K27").Value = 5198.1665
$ws.Range("M27").Value = -5091.1665
$ws.Range("H45").Value = 17000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("H55").Value = 101
$ws.Range("I55").Value = 87.85714
$ws.Range("J55").Value = 108.666664
$ws.Range("K55").Value = 87.85714
$ws.Range("L55").Value = 108.666664
$ws.Range("M55").Value = 85.14286
$ws.Range("N55").Value = -454.666664
$ws.Range("H82").Value = 2214
$ws.Range("I82").Value = 2080.875
$ws.Range("J82").Value = 2640
$ws.Range("K82").Value = 2080.875
$ws.Range("L82").Value = 2640
$ws.Range("M82").Value = -1719.875
$ws.Range("N82").Value = -3362
$ws.Range("H85").Value = 2214
$ws.Range("I85").Value = 2080.875
$ws.Range("J85").Value = 2640
$ws.Range("K85").Value = 2080.875
$ws.Range("L85").Value = 2640
$ws.Range("M85").Value = -832.875
$ws.Range("N85").Value = -5136
$ws.Range("H132").Value = 1508.0883
$ws.Range("I132").Value = 1223.1724
$ws.Range("K132").Value = 3669.5172
$ws.Range("M132").Value = -1139.5172
$ws.Range("M45").Value = $null

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1070.931
$ws.Range("I132").Value = 706.913
$ws.Range("J132").Value = 2466.3333
$ws.Range("K132").Value = 2120.739
$ws.Range("L132").Value = 7398.999899999999
$ws.Range("M132").Value = 409.261
$ws.Range("N132").Value = -12458.9999
